$d = $word.ActiveDocument

$newText = @'
Sara Donnelly Stobbr anne ry@@BR@@@@BR@@SKILLS@@BR@@Certified, Experienced Tutor, strong communication skilis ns STOKES 8nd Forme@@BR@@EXPERIENCE@@BR@@@@BR@@Loomis Gasin Dolphins, Del Ore High School Pool — swim instructor@@BR@@@@BR@@dune 6, 2022 - auly 10, 2022@@BR@@@@BR@@Taught kids ages 2-12 how to@@BR@@‘Tracked progress and wrote report@@BR@@ Aetea Se" Feguara to" ensure kids were cafe and not tn trouble while in the water@@BR@@@@BR@@Johnson Ranch Sports Club, Granite Bay, CA ~ Swim Instructor@@BR@@@@BR@@aUL strokes and water safety.@@BR@@@@BR@@Taught kids ages 2-12 how to@@BR@@@@BR@@Johnson Ranch Sports Club, Granite Bay, CA — swim Coach@@BR@@= Coach kids ages 4-18 in the sport of swimming.@@BR@@Sescn SL of the atrokes and advanced techniques when the awinners are ready for them.@@BR@@@@BR@@EDUCATION@@BR@@@@BR@@Del Ore High School, Loomis, CA — High School Diploma@@BR@@@@BR@@Sierra College, Rocklin, CA@@BR@@@@BR@@California state University, Sacramento, CA@@BR@@@@BR@@Santa Clara University, Santa Clara, CA@@BR@@@@BR@@Related coursewor!@@BR@@wpcogranming = Mechatrontes = Gateutus 242.3@@BR@@‘AWARDS@@BR@@‘2 Candidate for Valadictortan Academie Lerten@@BR@@2 fenos"Rete 2 SSaT*OE ToL TStacy@@BR@@@@BR@@ACTIVITIES@@BR@@
'@

# Remove the trailing newline that the here-string terminator introduces
$newText = $newText.TrimEnd("`r", "`n")

# Replace the placeholder token with an actual manual line break character (same as <w:br/>)
$newText = $newText.Replace("@@BR@@", [char]11)

$p = $d.Paragraphs.Item(2)
$p.Range.Text = $newText
